$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 413, shifting existing rows (and the table) down.
$ws.Rows.Item(413).Insert()

# Grow the Table1 definition (header + data range) to cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I437"))

# Populate the new row 413 with the new derived variable.
$ws.Range("A413").Value = "X12b"
$ws.Range("B413").Value = "der_vax_count"
$ws.Range("C413").Value = "Other"
$ws.Range("D413").Value = "Number of doses of vaccine received prior to COVID-19"
$ws.Range("E413").Value = "0 doses; 1 mrna dose; 1 non-mrna dose; 2 mrna doses; 2+ non-mrna doses; 3+ mrna doses; Other; Unknown (dose and/or timing)"
$ws.Range("F413").Value = "NA (missing)"
$ws.Range("G413").Value = "No"
$ws.Range("H413").Value = "No"
$ws.Range("I413").Value = ""

$ws.Rows.Item(413).RowHeight = 46

# Two other existing rows grow taller to fit wrapped text.
$ws.Rows.Item(36).RowHeight = 46
$ws.Rows.Item(38).RowHeight = 46
